$d = $word.ActiveDocument

$replacements = @{
    "1. Introduction"          = "1. 1. Introduction"
    "2. Data Cleaning"         = "2. 2. Data Cleaning"
    "3. EDA"                   = "3. 3. EDA"
    "4. Skill Gap Analysis"    = "4. 4. Skill Gap Analysis"
}

foreach ($p in $d.Paragraphs) {
    $paraText = $p.Range.Text.TrimEnd("`r", "`a")
    if ($replacements.ContainsKey($paraText)) {
        $p.Range.Text = $replacements[$paraText]
    }
}
